$props = $ppt | Get-Member
Write-Output $props
